# Update cryptos list values (price/volume%) for rows 2-51 to match the latest
# scrape, and fix the ordering swap for Hedera/TrustWalletToken (rows 38-39)
# and EnergySwap/Decentraland (rows 47-48).
#
# Price cells that look like plain decimals (e.g. "1.001") are forced to Text
# format before assignment so Excel keeps them as text instead of coercing them
# into numbers (matching the original inlineStr/text cell type).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.625.15"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.850.50"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.48"
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.63"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07319"
$ws.Range("E10").Value = "  -2.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8883"
$ws.Range("E11").Value = "  -5.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.84"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "1.917.30"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.366"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.559"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06912"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.11"
$ws.Range("E18").Value = "  -2.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008908"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.49"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "27.648.93"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.994"
$ws.Range("E23").Value = "  -2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.68"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").Value = "2.088.91"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.969"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.75"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.02"
$ws.Range("E28").Value = "  +2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "122.23"
$ws.Range("E29").Value = "  +7.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.266"
$ws.Range("E30").Value = "  -5.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.924"
$ws.Range("E31").Value = "  +12.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08944"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7671"
$ws.Range("E33").Value = "  -6.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.589"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.990"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.105"
$ws.Range("E36").Value = "  -6.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05396"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.096"
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01952"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.810"
$ws.Range("E41").Value = "  -4.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.935"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5124"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1659"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.306"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06583"
$ws.Range("E46").Value = "  -2.48%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4776"
$ws.Range("E47").Value = "  -2.46%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.40"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.54"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9999"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.634"
$ws.Range("E51").Value = "  -2.43%  "
